$d = $word.ActiveDocument

# The paragraph currently holds three runs: "<id>", "p167r_1", "</id>".
# Collapse them into a single run (picking up the surrounding Courier-New
# formatting) containing the combined text "<id>p167r_1</id>".
$rng = $d.Content
$null = $rng.Find.Execute("<id>p167r_1</id>", $true, $false, $false, $false, `
                           $false, $true, 1, $false, "<id>p167r_1</id>", 2)
